$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties in columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style of the existing header cells (e.g. AC1) for the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill team record (Wins=96, Losses=66, Ties=0) for every data row (2..42)
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 96
    $ws.Cells.Item($r, 31).Value = 66
    $ws.Cells.Item($r, 32).Value = 0
}
